# StatsBasePerLevel.xlsx edit:
#  - add an "Item" (shield index) column right after the shield (H) column
#  - change the way the shield (H column) base stat is calculated
#  - rename abbreviated stat headers to full words

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the abbreviated header cells to their full word equivalents.
$ws.Range("D1").Value() = "endurance"
$ws.Range("E1").Value() = "power"
$ws.Range("F1").Value() = "chance"
$ws.Range("G1").Value() = "velocity"

# Insert a new column right after H (shield) to hold the new "shield index"
# item column. This shifts the old I (EXP_Total) and J (Exp(n)-Exp(n-1))
# columns one position to the right (J and K).
$ws.Columns.Item(9).Insert()

# New header for the inserted "shield index" column (I).
$ws.Range("I1").Value() = "shield"

# Change the way the shield base stat (column H) is calculated, and fill
# in the new "shield index" column (I) based on it. Row 2 is written
# separately from rows 3:31 so the resulting formula layout mirrors how it
# was originally entered (a single cell, then the rest filled down).
$ws.Range("H2").Formula = "=(A2+8)*(A2-1)/2*5+200"
$ws.Range("H3:H31").Formula = "=(A3+8)*(A3-1)/2*5+200"

$ws.Range("I2").Formula = "=H2/200"
$ws.Range("I3:I31").Formula = "=H3/200"

# Re-fill the "Exp(n)-Exp(n-1)" formulas (now in column K) below the first
# two rows, since they were shifted out of column J by the column insert.
$ws.Range("K4:K31").Formula = "=J4-J3"

# Update the selection to match the edited range.
$ws.Range("H2:H31").Select()
